# Applies the "Add descriptions titles" commit:
#  - Metadata sheet: fill in Title + Description values, bump Date timestamp
#  - Elements sheet: give the root "Extension" row a real Short/Definition
#    (replacing the generic boilerplate text) and clear its RIM mapping

$wb = $excel.ActiveWorkbook

$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B5").Value = "DMI Facture"
$meta.Range("B8").Value = "2026-02-25T08:15:31+00:00"
$meta.Range("B12").Value = "Extension créée dans ce volet pour représenter une facture."

$elements = $wb.Worksheets.Item("Elements")
$elements.Range("L2").Value = "DMI Facture"
$elements.Range("M2").Value = "Extension créée dans ce volet pour représenter une facture."
$elements.Range("AK2").Value = ""
